$p = $ppt.ActivePresentation

# slide2.xml -> Slides.Item(2)
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 1 – Titre" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "L’application de course compétitive en temps réel" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Sous‑titre :" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " Transformez chaque sortie running en duel excitant." }
    }
}

# slide3.xml -> Slides.Item(3)
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 2 – Problème" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Courir seul, une motivation en berne" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "De nombreux coureurs se sentent démotivés seuls" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Manque de défi et d’interaction" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "Difficulté à mesurer ses progrès face à d’autres" }
        if ($pIdx -eq 5 -and $rIdx -eq 1) { $run.Text = "Risque d’abandon des objectifs sportifs" }
        if ($pIdx -eq 6 -and $rIdx -eq 1) { $run.Text = "➡ Résultat : monotonie, perte de motivation, arrêt de la pratique" }
    }
}

# slide4.xml -> Slides.Item(4)
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 3 – Nos différenciateurs" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Ce qui rend GeoRace unique" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Duel en temps réel vers un point d’arrivée équidistant" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " → équité garantie entre coureurs" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Matchmaking ELO" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " → adversaires de niveau comparable, progression visible" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "GPS en direct & suivi synchronisé" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " → expérience compétitive immédiate" }
        if ($pIdx -eq 5 -and $rIdx -eq 1) { $run.Text = "Confidentialité et sécurité" }
        if ($pIdx -eq 5 -and $rIdx -eq 2) { $run.Text = " → partage de position limité aux duels actifs" }
        if ($pIdx -eq 6 -and $rIdx -eq 1) { $run.Text = "Anti‑triche intégré" }
        if ($pIdx -eq 6 -and $rIdx -eq 2) { $run.Text = " (détection d’anomalies GPS, vérifications de parcours)" }
        if ($pIdx -eq 7 -and $rIdx -eq 1) { $run.Text = "Fonctions sociales" }
        if ($pIdx -eq 7 -and $rIdx -eq 2) { $run.Text = " → profils, chat, badges, classements locaux" }
    }
}

# slide5.xml -> Slides.Item(5)
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 4 – Pourquoi les utilisateurs nous rejoindront" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Les raisons de passer à GeoRace" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Expérience de duel unique et équitable" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " = plus de fun et de challenge" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Progression visible (ELO)" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " = motivation à revenir" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "Événements locaux & tournois" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " = opportunités sociales et récompenses" }
        if ($pIdx -eq 5 -and $rIdx -eq 1) { $run.Text = "Transfert facile des amis" }
        if ($pIdx -eq 5 -and $rIdx -eq 2) { $run.Text = " via invitations et duels immédiats" }
    }
}

# slide6.xml -> Slides.Item(6)
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 5 – Stratégie d’acquisition (comment voler des utilisateurs aux concurrents)" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Tactiques concrètes pour attirer les utilisateurs" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Ciblage local" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " : organiser micro‑événements et challenges dans les parcs et campus" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Partenariats" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " avec clubs de course, coachs, influenceurs locaux" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "Intégrations & porte d’entrée" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " : import simplifié depuis Strava/Runkeeper + import d’amis" }
        if ($pIdx -eq 5 -and $rIdx -eq 1) { $run.Text = "Programmes d’incitation" }
        if ($pIdx -eq 5 -and $rIdx -eq 2) { $run.Text = " : duels gratuits/bonus initial ELO, récompenses, badges exclusifs" }
        if ($pIdx -eq 6 -and $rIdx -eq 1) { $run.Text = "Campagnes social media" }
        if ($pIdx -eq 6 -and $rIdx -eq 2) { $run.Text = " montrant duels en direct + récits d’utilisateurs" }
        if ($pIdx -eq 7 -and $rIdx -eq 1) { $run.Text = "Fonctionnalités différenciantes en onboarding" }
        if ($pIdx -eq 7 -and $rIdx -eq 2) { $run.Text = " (démonstration duel en 30s)" }
    }
}

# slide7.xml -> Slides.Item(7)
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 6 – Concept innovant : Le mode duel" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Duels en temps réel pour tous" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = "DUELS EN TEMPS RÉEL" }
        if ($pIdx -eq 2 -and $rIdx -eq 3) { $run.Text = " Affrontez d’autres coureurs à proximité dans des courses spontanées" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = "POINT D’ARRIVÉE ÉQUITABLE" }
        if ($pIdx -eq 3 -and $rIdx -eq 3) { $run.Text = " Chaque duel se termine à un point équidistant pour tous les participants" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = "COMPÉTITION JUSTE" }
        if ($pIdx -eq 4 -and $rIdx -eq 3) { $run.Text = " Système ELO garantit des matchs équilibrés quel que soit votre niveau" }
    }
}

# slide8.xml -> Slides.Item(8)
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 7 – Fonctionnement d’un duel" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "3 étapes simples" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "01 - DÉTECTION DES COUREURS" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " Identification automatique des coureurs actifs à proximité" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "02 - POINT D’ARRIVÉE ÉQUITABLE" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " Calcul algorithmique d’un point accessible à distance égale pour tous" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "03 - DÉPART SIMULTANÉ" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " Compte à rebours synchronisé après acceptation mutuelle" }
    }
}

# slide9.xml -> Slides.Item(9)
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 8 – Un duel GeoRace" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Scénario typique :" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Vous commencez votre course" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " → l’app détecte 3 coureurs dans un rayon de 800m" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Proposition de duel" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " → notification “Paul (ELO 1420) vous défie”" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "Vous acceptez" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " → le système calcule un point d’arrivée à 1,2km pour chacun" }
        if ($pIdx -eq 5 -and $rIdx -eq 1) { $run.Text = "Départ dans 10 secondes" }
        if ($pIdx -eq 5 -and $rIdx -eq 2) { $run.Text = " → compte à rebours" }
        if ($pIdx -eq 6 -and $rIdx -eq 1) { $run.Text = "Course en direct" }
        if ($pIdx -eq 6 -and $rIdx -eq 2) { $run.Text = " → GPS vous guide, position de l’adversaire visible" }
        if ($pIdx -eq 7 -and $rIdx -eq 1) { $run.Text = "Arrivée" }
        if ($pIdx -eq 7 -and $rIdx -eq 2) { $run.Text = " → +25 ELO, badge “5 victoires consécutives” débloqué" }
    }
}

# slide10.xml -> Slides.Item(10)
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 9 – Démo rapide (points clés pour la présentation)" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Script de présentation - Points à appuyer" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Point d’arrivée équidistant" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " : “Chaque coureur parcourt exactement la même distance - l’équité totale”" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Matchmaking ELO" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " : “Vous affrontez toujours quelqu’un de votre niveau - ni trop facile, ni impossible”" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "Duel en temps réel" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " : “Pas de comparaison après coup - c’est une vraie course, maintenant”" }
        if ($pIdx -eq 5 -and $rIdx -eq 1) { $run.Text = "Sécurité et confidentialité" }
        if ($pIdx -eq 5 -and $rIdx -eq 2) { $run.Text = " : “Votre position n’est visible que pendant les duels actifs”" }
        if ($pIdx -eq 6 -and $rIdx -eq 1) { $run.Text = "Anti-triche" }
        if ($pIdx -eq 6 -and $rIdx -eq 2) { $run.Text = " : “Détection automatique des anomalies GPS - pas de tricheurs”" }
        if ($pIdx -eq 7 -and $rIdx -eq 1) { $run.Text = "Engagement social" }
        if ($pIdx -eq 7 -and $rIdx -eq 2) { $run.Text = " : “Profils, chat, badges, tournois locaux - une vraie communauté”" }
        if ($pIdx -eq 8 -and $rIdx -eq 1) { $run.Text = "Message clé" }
        if ($pIdx -eq 8 -and $rIdx -eq 2) { $run.Text = " : Le duel équitable en temps réel change tout - ce n’est plus du tracking, c’est de la compétition !" }
    }
}

# slide11.xml -> Slides.Item(11)
$s = $p.Slides.Item(11)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 10 – Aspect social" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Construisez une communauté de coureurs" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "01 - MATCHMAKING ÉQUILIBRÉ" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " Connexion automatique entre coureurs de niveaux similaires grâce à l’ELO" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "02 - DISCUSSION ET PARTAGE" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " Chat en direct, ajout d’amis, profils détaillés avec statistiques" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "03 - ÉVÉNEMENTS COMMUNAUTAIRES" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " Courses organisées, défis collectifs, ligues locales par ville/quartier" }
    }
}

# slide12.xml -> Slides.Item(12)
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 11 – Système de classement ELO" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Compétition équitable et progression motivante" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "01 - MATCHS ÉQUILIBRÉS" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " Algorithme ELO vous oppose à des coureurs de niveau proche" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "02 - PROGRESSION VISIBLE" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " Suivez votre montée dans les classements (Bronze → Argent → Or → Platine)" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "03 - RÉCOMPENSES ET BADGES" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " Débloquez des achievements à chaque niveau atteint" }
    }
}

# slide13.xml -> Slides.Item(13)
$s = $p.Slides.Item(13)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 12 – Bénéfices de GeoRace" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Plus qu’une app de running" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "01 - MOTIVATION ACCRUE" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " Chaque sortie devient un défi concret et excitant" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "02 - AMÉLIORATION DES PERFORMANCES" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " La compétition pousse naturellement au dépassement de soi" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "03 - BIEN-ÊTRE SOCIAL" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " Rencontrez et interagissez avec d’autres coureurs passionnés" }
    }
}

# slide14.xml -> Slides.Item(14)
$s = $p.Slides.Item(14)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 13 – Fonctionnalités clés" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Technologie au service de la performance" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "GPS haute précision" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " : tracking en temps réel pendant les duels" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Algorithme intelligent" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " : calcul de points équidistants accessibles" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "Sécurité" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " : position visible uniquement en duel actif" }
        if ($pIdx -eq 5 -and $rIdx -eq 1) { $run.Text = "Anti-triche" }
        if ($pIdx -eq 5 -and $rIdx -eq 2) { $run.Text = " : détection de vitesse anormale et validation communautaire" }
        if ($pIdx -eq 6 -and $rIdx -eq 1) { $run.Text = "Optimisation batterie" }
        if ($pIdx -eq 6 -and $rIdx -eq 2) { $run.Text = " : mode économie d’énergie" }
    }
}

# slide15.xml -> Slides.Item(15)
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 14 – Cible" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Qui utilise GeoRace ?" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = "Coureurs réguliers" }
        if ($pIdx -eq 2 -and $rIdx -eq 3) { $run.Text = " : cherchent à améliorer leurs performances" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = "Coureurs occasionnels" }
        if ($pIdx -eq 3 -and $rIdx -eq 3) { $run.Text = " : ont besoin de motivation pour rester réguliers" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = "Coureurs compétitifs" }
        if ($pIdx -eq 4 -and $rIdx -eq 3) { $run.Text = " : veulent tester leur niveau entre les courses officielles" }
        if ($pIdx -eq 5 -and $rIdx -eq 2) { $run.Text = "Nouveaux coureurs" }
        if ($pIdx -eq 5 -and $rIdx -eq 3) { $run.Text = " : cherchent une communauté motivante" }
    }
}

# slide16.xml -> Slides.Item(16)
$s = $p.Slides.Item(16)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 15 – Différenciation vs concurrents" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Ce que les autres apps ne font pas :" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " → réseau social, pas de compétition temps réel" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " → tracking + défis asynchrones" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " → GPS basique sans interaction" }
        if ($pIdx -eq 5 -and $rIdx -eq 1) { $run.Text = "GeoRace combine :" }
        if ($pIdx -eq 6 -and $rIdx -eq 1) { $run.Text = "Duels spontanés en temps réel + Point d’arrivée équitable + Matchmaking ELO" }
        if ($pIdx -eq 7 -and $rIdx -eq 1) { $run.Text = "= Concept unique sur le marché" }
    }
}

# slide17.xml -> Slides.Item(17)
$s = $p.Slides.Item(17)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 16 – Modèle économique" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Freemium avec valeur ajoutée" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Gratuit :" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " * 3 duels par jour * Classement ELO * Statistiques de base * Événements communautaires" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Premium (4,99€/mois) :" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " * Duels illimités * Statistiques avancées * Personnalisation profil * Mode entraînement contre vos records * Pas de publicités" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "Revenus additionnels :" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " * Partenariats marques (Nike, Adidas) * Événements premium avec prix" }
    }
}

# slide18.xml -> Slides.Item(18)
$s = $p.Slides.Item(18)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Déploiement progressif" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Phase 1 (6 mois) - MVP" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " * Duels 1v1 en temps réel * Système ELO * Test dans 2-3 villes pilotes * " }
        if ($pIdx -eq 2 -and $rIdx -eq 3) { $run.Text = "Objectif : 500 utilisateurs actifs" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Phase 2 (12 mois) - Gamification" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " * Badges, statistiques, profils * Événements communautaires * " }
        if ($pIdx -eq 3 -and $rIdx -eq 3) { $run.Text = "Objectif : 5 000 utilisateurs" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "Phase 3 (24 mois) - Scale" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " * Version Premium * Duels multi-joueurs * Expansion nationale * " }
        if ($pIdx -eq 4 -and $rIdx -eq 3) { $run.Text = "Objectif : 50 000 utilisateurs" }
    }
}

# slide19.xml -> Slides.Item(19)
$s = $p.Slides.Item(19)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 18 – Défis & Solutions" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "Anticiper les obstacles" }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Défi : Sécurité des coureurs" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " → Alertes de sécurité, mode “course prudente”" }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Défi : Zones rurales (peu d’utilisateurs)" }
        if ($pIdx -eq 3 -and $rIdx -eq 2) { $run.Text = " → Mode asynchrone contre “ghost runners”" }
        if ($pIdx -eq 4 -and $rIdx -eq 1) { $run.Text = "Défi : Triche (vélo, voiture)" }
        if ($pIdx -eq 4 -and $rIdx -eq 2) { $run.Text = " → Détection vitesse anormale + validation communautaire" }
        if ($pIdx -eq 5 -and $rIdx -eq 1) { $run.Text = "Défi : Vie privée" }
        if ($pIdx -eq 5 -and $rIdx -eq 2) { $run.Text = " → Position visible uniquement en duel, blocage d’utilisateurs" }
    }
}

# slide20.xml -> Slides.Item(20)
$s = $p.Slides.Item(20)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "GeoRace ne se contente pas de tracker vos courses. Nous transformons la course à pied en une " }
        if ($pIdx -eq 1 -and $rIdx -eq 2) { $run.Text = "expérience sociale, compétitive et motivante" }
        if ($pIdx -eq 1 -and $rIdx -eq 3) { $run.Text = " accessible à tous, partout, à tout moment." }
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Notre ambition :" }
        if ($pIdx -eq 2 -and $rIdx -eq 2) { $run.Text = " Devenir la référence de la course compétitive en temps réel et créer une communauté mondiale de coureurs passionnés par le défi." }
    }
}

# slide21.xml -> Slides.Item(21)
$s = $p.Slides.Item(21)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 1 -and $rIdx -eq 1) { $run.Text = "🟦 Slide 20 – Phrase finale" }
    }
}
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$pIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $pIdx++
    $rIdx = 0
    foreach ($run in $para.Runs()) {
        $rIdx++
        if ($pIdx -eq 2 -and $rIdx -eq 1) { $run.Text = "Chaque course est un duel. Chaque duel est une opportunité. Transformez votre motivation." }
        if ($pIdx -eq 3 -and $rIdx -eq 1) { $run.Text = "Prêt à relever le défi ?" }
    }
}
